$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold, border, centered) used by the existing header
# cells in row 1 onto the two new header cells before setting their text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells (row 2), using the default/unstyled format like H2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
